$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2: "David Darquea" -> "David"
$ws.Range("A2").Value = "David"

# Update C2: 17 -> 20
$ws.Range("C2").Value = 20

# Delete row 3 (Marcos Rodriguez) entirely
$ws.Rows(3).Delete()
